$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.608.23"
$ws.Range("E2").Value = "  -2.22%  "

# Row 3
$ws.Range("D3").Value = "1.677.23"
$ws.Range("E3").Value = "  -1.55%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.54%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.82"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5212"
$ws.Range("E6").Value = "  -1.55%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +0.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06538"
$ws.Range("E8").Value = "  -0.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2587"
$ws.Range("E9").Value = "  -2.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.23"
$ws.Range("E10").Value = "  -2.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07700"
$ws.Range("E11").Value = "  +0.85%  "

# Row 12
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.914.20"
$ws.Range("E12").Value = "  -1.23%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.686.91"
$ws.Range("E13").Value = "  -1.39%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.354"
$ws.Range("E14").Value = "  -4.66%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5626"
$ws.Range("E15").Value = "  -1.54%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0@8056"
$ws.Range("D16").Replace("@", [string]([char]0x2085)) | Out-Null
$ws.Range("E16").Value = "  -1.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.60"
$ws.Range("E17").Value = "  -2.57%  "

# Row 18
$ws.Range("D18").Value = "26.706.40"
$ws.Range("E18").Value = "  -1.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.55"
$ws.Range("E19").Value = "  -0.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.008"
$ws.Range("E20").Value = "  +0.40%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.515"
$ws.Range("E21").Value = "  -2.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.18"
$ws.Range("E22").Value = "  -2.24%  "

# Row 23
$ws.Range("E23").Value = "  -0.30%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.009"
$ws.Range("E24").Value = "  +0.58%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.69"
$ws.Range("E25").Value = "  +1.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.738"
$ws.Range("E26").Value = "  -0.44%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1174"
$ws.Range("E27").Value = "  -3.46%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.080"
$ws.Range("E28").Value = "  -2.20%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.90"
$ws.Range("E29").Value = "  -2.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05264"
$ws.Range("E30").Value = "  -1.63%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.274"
$ws.Range("E31").Value = "  -0.98%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.384"
$ws.Range("E32").Value = "  -3.18%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.260"
$ws.Range("E33").Value = "  -4.21%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.605"
$ws.Range("E34").Value = "  -1.56%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.784"
$ws.Range("E35").Value = "  -3.06%  "

# Row 36
$ws.Range("E36").Value = "  -1.07%  "

# Row 37
$ws.Range("E37").Value = "  -1.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5766"
$ws.Range("E38").Value = "  -1.28%  "

# Row 39
$ws.Range("D39").Value = "1.165.93"
$ws.Range("E39").Value = "  +12.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01617"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.726"
$ws.Range("E42").Value = "  -2.21%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8328"
$ws.Range("E43").Value = "  -0.56%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.87"
$ws.Range("E44").Value = "  -0.95%  "

# Row 45
$ws.Range("D45").Value = "1.821.77"
$ws.Range("E45").Value = "  -1.31%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0@111"
$ws.Range("D46").Replace("@", [string]([char]0x2088)) | Out-Null
$ws.Range("E46").Value = "  -3.79%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.25"
$ws.Range("E47").Value = "  -2.75%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4501"
$ws.Range("E48").Value = "  +0.24%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  +0.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.994"
$ws.Range("E50").Value = "  -0.99%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05180"
$ws.Range("E51").Value = "  -1.13%  "
